$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data row (row 5) with the failed-login test fixture values
# (stored base64-encoded, matching the existing rows' convention).
# Value2 is used so the all-digit password string isn't coerced to a number.
$ws.Cells.Item(5, 1).Value2 = "dXN1YXJpb0ZhbGxpZG8="
$ws.Cells.Item(5, 2).Value2 = "MTIzNA=="

# Update the active selection to match the new last cell
$ws.Range("B5").Select()
